$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The MSquared macro gained a new NET parameter, so the matching test rows'
# Description text (column B) is updated to document it. Look each test row
# up by its name in column A so the edit does not depend on an assumed row
# number.
$updates = @{
    "Msquared1" = "Test Msquared with method=discrete, scale=252, NET=TRUE"
    "Msquared2" = "Test Msquared with method=discrete, scale=1,NET=FALSE"
    "Msquared3" = "Test Msquared with method=log, scale=4,NET=TRUE"
    "Msquared4" = "Test Msquared with method=log, scale=12,NET=FALSE"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

foreach ($name in $updates.Keys) {
    $found = $ws.Range("A1:A$lastRow").Find($name)
    if ($found -ne $null) {
        $row = $found.Row
        $ws.Cells.Item($row, 2).Value = $updates[$name]
    }
}

# Reflect the cell that was being worked on / scrolled to when the file was
# last saved.
$ws.Range("B35").Select()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
